$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Value
    )
    $rng = $ws.Range($CellRef)
    # Force text storage so numeric-looking strings (e.g. "245.32") are not
    # auto-converted to numbers by Excel's input parsing, then restore the
    # default "Normal" style so we don't leave a stray number-format style
    # behind on the cell (matches the original inline-string cells, which
    # carry no explicit style).
    $rng.NumberFormat = "@"
    $rng.Value = $Value
    $rng.Style = "Normal"
}

# Column D ("Price") updates
Set-TextValue "D2"  "245.32"
Set-TextValue "D3"  "23.92"
Set-TextValue "D4"  "5.204"
Set-TextValue "D5"  "0.05743"
Set-TextValue "D6"  "6.455"
Set-TextValue "D7"  "3.213"
Set-TextValue "D8"  "0.8143"
Set-TextValue "D9"  "0.8674"
Set-TextValue "D11" "0.06996"
Set-TextValue "D13" "0.03017"
Set-TextValue "D14" "0.09330"
Set-TextValue "D15" "3.823"
Set-TextValue "D16" "0.001528"
Set-TextValue "D17" "0.04728"
Set-TextValue "D18" "0.0006008"
Set-TextValue "D19" "0.006210"
Set-TextValue "D20" "0.001236"
Set-TextValue "D22" "0.00008697"
Set-TextValue "D24" "2.153"
Set-TextValue "D26" "0.1330"
Set-TextValue "D27" "0.0002328"
Set-TextValue "D40" "0.03717"
Set-TextValue "D41" "0.006218"
Set-TextValue "D43" "0.002299"
Set-TextValue "D44" "0.007458"
Set-TextValue "D45" "0.00005376"
Set-TextValue "D46" "0.00000000750"
Set-TextValue "D47" "0.4398"
Set-TextValue "D48" "0.002043"
Set-TextValue "D49" "0.00002099"
Set-TextValue "D50" "0.0001999"

# Column E ("Volume(1h)") updates - plain alphanumeric text, no numeric
# coercion risk, so a direct assignment is sufficient.
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"
$ws.Range("E48").Value = "47BOLOBOLO"

Write-Output "Applied symbol list update."
